$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.514.62"
$ws.Range("E2").Value = "  +0.16%  "
$ws.Range("D3").Value = "1.822.52"
$ws.Range("E3").Value = "  -0.08%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'315.27"
$ws.Range("E5").Value = "  -0.53%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("D7").Value = "'0.5103"
$ws.Range("E7").Value = "  -5.36%  "
$ws.Range("D8").Value = "'0.3945"
$ws.Range("E8").Value = "  -2.30%  "
$ws.Range("D9").Value = "'0.08183"
$ws.Range("E9").Value = "  +6.60%  "
$ws.Range("E10").Value = "  -1.21%  "
$ws.Range("D11").Value = "'41.66"
$ws.Range("E11").Value = "  -0.57%  "
$ws.Range("D12").Value = "'6.335"
$ws.Range("E12").Value = "  +0.06%  "
$ws.Range("D13").Value = "'21.10"
$ws.Range("E13").Value = "  +0.85%  "
$ws.Range("E14").Value = "  -0.02%  "
$ws.Range("D15").Value = "'7.513"
$ws.Range("E15").Value = "  -1.72%  "
$ws.Range("D16").Value = "1.822.17"
$ws.Range("E16").Value = "  -0.09%  "
$ws.Range("D17").Value = "'0.00001132"
$ws.Range("E17").Value = "  +3.31%  "
$ws.Range("D18").Value = "'92.56"
$ws.Range("E18").Value = "  +3.03%  "
$ws.Range("D19").Value = "'0.06658"
$ws.Range("E19").Value = "  +0.75%  "
$ws.Range("D20").Value = "'17.80"
$ws.Range("E20").Value = "  +0.57%  "
$ws.Range("D21").Value = "'1.000"
$ws.Range("E21").Value = "  -0.08%  "
$ws.Range("D22").Value = "'6.093"
$ws.Range("E22").Value = "  +0.42%  "
$ws.Range("D23").Value = "28.538.11"
$ws.Range("E23").Value = "  +0.22%  "
$ws.Range("E24").Value = "  +2.23%  "
$ws.Range("E25").Value = "  -0.21%  "
$ws.Range("D26").Value = "'21.37"
$ws.Range("E26").Value = "  +3.03%  "
$ws.Range("D27").Value = "'156.53"
$ws.Range("E27").Value = "  -0.81%  "
$ws.Range("D28").Value = "2.030.65"
$ws.Range("E28").Value = "  -0.20%  "
$ws.Range("D29").Value = "'2.404"
$ws.Range("E29").Value = "  -2.28%  "
$ws.Range("D30").Value = "'125.97"
$ws.Range("E30").Value = "  +1.65%  "
$ws.Range("E31").Value = "  -0.64%  "
$ws.Range("D32").Value = "'0.1094"
$ws.Range("E32").Value = "  -1.23%  "
$ws.Range("D33").Value = "'5.766"
$ws.Range("E33").Value = "  +1.58%  "
$ws.Range("D34").Value = "'3.655"
$ws.Range("E34").Value = "  +0.38%  "
$ws.Range("D35").Value = "'0.07070"
$ws.Range("E35").Value = "  -3.66%  "
$ws.Range("D36").Value = "'0.2228"
$ws.Range("E36").Value = "  -0.55%  "
$ws.Range("D37").Value = "'0.02356"
$ws.Range("E37").Value = "  +0.53%  "
$ws.Range("D38").Value = "'5.261"
$ws.Range("E38").Value = "  +0.97%  "
$ws.Range("D39").Value = "'8.841"
$ws.Range("E39").Value = "  -0.22%  "
$ws.Range("D40").Value = "'0.6314"
$ws.Range("E40").Value = "  +0.31%  "
$ws.Range("E41").Value = "  -0.49%  "
$ws.Range("D42").Value = "'1.181"
$ws.Range("E42").Value = "  -0.64%  "
$ws.Range("D43").Value = "'1.000"
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("D44").Value = "'1.398"
$ws.Range("E44").Value = "  -0.15%  "
$ws.Range("D45").Value = "'13.48"
$ws.Range("E45").Value = "  +0.29%  "
$ws.Range("D46").Value = "'0.5931"
$ws.Range("E46").Value = "  +1.36%  "
$ws.Range("D47").Value = "'3.732"
$ws.Range("E47").Value = "  +0.89%  "
$ws.Range("E48").Value = "  -0.25%  "
$ws.Range("D49").Value = "'1.987"
$ws.Range("E49").Value = "  -0.79%  "
$ws.Range("D50").Value = "'1.186"
$ws.Range("E50").Value = "  -1.08%  "
$ws.Range("D51").Value = "'0.06896"
$ws.Range("E51").Value = "  +0.27%  "
